# Update "想去人数" (column F) counts on the sheets that list event details.
# The workbook has 4 sheets: 展览 (1), 演出 (2), 本地生活 (3), 全部类型 (4).
# Sheets 1 and 4 share the same event rows / column layout and both need
# the same column-F updates, mirroring the diff.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 185
    3  = 227
    5  = 785
    6  = 233
    7  = 5778
    8  = 26
    10 = 97
    11 = 40
    14 = 178
    15 = 326
    16 = 24
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
